$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "41.457.55"
$ws.Cells.Item(2, 5).Value = "  +0.61%  "

# Row 3
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "2.203.02"
$ws.Cells.Item(3, 5).Value = "  -0.77%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  +0.00%  "

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "256.10"
$ws.Cells.Item(5, 5).Value = "  +4.26%  "

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "0.630"
$ws.Cells.Item(6, 5).Value = "  +0.28%  "

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "69.29"
$ws.Cells.Item(7, 5).Value = "  -1.39%  "

# Row 8
$ws.Cells.Item(8, 5).Value = "  -0.02%  "

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.579"
$ws.Cells.Item(9, 5).Value = "  +4.06%  "

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "38.01"
$ws.Cells.Item(10, 5).Value = "  -2.40%  "

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "58.83"
$ws.Cells.Item(11, 5).Value = "  +0.83%  "

# Row 12
$ws.Cells.Item(12, 5).Value = "  -0.24%  "

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "7.15"
$ws.Cells.Item(13, 5).Value = "  +5.80%  "

# Row 14
$ws.Cells.Item(14, 5).Value = "  +0.20%  "

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "2.528.79"
$ws.Cells.Item(15, 5).Value = "  -0.73%  "

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "0.884"
$ws.Cells.Item(16, 5).Value = "  +4.80%  "

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "14.62"
$ws.Cells.Item(17, 5).Value = "  -1.44%  "

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "2.181.91"
$ws.Cells.Item(18, 5).Value = "  -1.82%  "

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "41.459.00"
$ws.Cells.Item(19, 5).Value = "  +0.46%  "

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "0.0₃0966"
$ws.Cells.Item(20, 5).Value = "  +1.33%  "

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "6.23"
$ws.Cells.Item(21, 5).Value = "  +2.53%  "

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "72.63"
$ws.Cells.Item(22, 5).Value = "  +0.17%  "

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "234.82"
$ws.Cells.Item(23, 5).Value = "  +1.24%  "

# Row 24
$ws.Cells.Item(24, 5).Value = "  -1.97%  "

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "11.93"
$ws.Cells.Item(25, 5).Value = "  +20.89%  "

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "3.90"
$ws.Cells.Item(26, 5).Value = "  +6.10%  "

# Row 27
$ws.Cells.Item(27, 5).Value = "  +0.15%  "

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "2.52"
$ws.Cells.Item(28, 5).Value = "  +3.67%  "

# Row 29
$ws.Cells.Item(29, 2).Value = "Toncoin"
$ws.Cells.Item(29, 3).Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "2.18"
$ws.Cells.Item(29, 5).Value = "  -0.14%  "

# Row 30
$ws.Cells.Item(30, 2).Value = "Monero"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "169.84"
$ws.Cells.Item(30, 5).Value = "  -1.40%  "

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "20.77"
$ws.Cells.Item(31, 5).Value = "  +1.14%  "

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "0.118"
$ws.Cells.Item(32, 5).Value = "  -0.69%  "

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "0.0759"
$ws.Cells.Item(33, 5).Value = "  +6.55%  "

# Row 34
$ws.Cells.Item(34, 5).Value = "  -0.03%  "

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "5.50"
$ws.Cells.Item(35, 5).Value = "  +4.65%  "

# Row 36
$ws.Cells.Item(36, 2).Value = "RenderToken"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "4.19"
$ws.Cells.Item(36, 5).Value = "  +7.21%  "

# Row 37
$ws.Cells.Item(37, 2).Value = "InjectiveProtocol"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "26.18"
$ws.Cells.Item(37, 5).Value = "  +7.74%  "

# Row 38
$ws.Cells.Item(38, 2).Value = "Filecoin"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "4.65"
$ws.Cells.Item(38, 5).Value = "  +0.88%  "

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.0302"
$ws.Cells.Item(39, 5).Value = "  +9.51%  "

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "2.23"
$ws.Cells.Item(40, 5).Value = "  -1.30%  "

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "12.28"
$ws.Cells.Item(41, 5).Value = "  +12.52%  "

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "5.73"
$ws.Cells.Item(42, 5).Value = "  -2.20%  "

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "64.14"
$ws.Cells.Item(43, 5).Value = "  -2.36%  "

# Row 44
$ws.Cells.Item(44, 5).Value = "  -1.00%  "

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "0.201"
$ws.Cells.Item(45, 5).Value = "  -1.94%  "

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "1.24"
$ws.Cells.Item(46, 5).Value = "  +13.38%  "

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "8.74"
$ws.Cells.Item(47, 5).Value = "  -1.08%  "

# Row 48
$ws.Cells.Item(48, 5).Value = "  +2.13%  "

# Row 49
$ws.Cells.Item(49, 5).Value = "  +0.39%  "

# Row 50
$ws.Cells.Item(50, 5).Value = "  +1.22%  "

# Row 51
$ws.Cells.Item(51, 2).Value = "NEARProtocol"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "2.30"
$ws.Cells.Item(51, 5).Value = "  +1.99%  "
